# The workbook stores file paths (image/sound asset references used by the
# experiment) as forward-slash paths, e.g. "Condition/CS+2.BMP",
# "Sound/silent.wav". The author reorganised these into a local
# "original huang experiment" folder on Windows, so the same references now
# use Windows-style backslash separators, e.g. "Condition\CS+2.BMP",
# "Sound\silent.wav". Apply that rename across every cell that references
# these asset paths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("Condition/", "Condition\")
$ws.Cells.Replace("Sound/", "Sound\")
